$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 15 with the new "token precision vs decimals" comparison notes
# (write F before E so the shared-string table picks up the same ordering
# as the authored workbook)
$ws.Range("A15").Value = "name of decimals in a token amount"
$ws.Range("B15").Value = "token precision (4 max.)"
$ws.Range("C15").Value = "token decimals (18 max.)"
$ws.Range("F15").Value = "lowest value, unit = 10^(-18), wei"
$ws.Range("E15").Value = "lowest value, unit = 0.0001, -"

# Row 15 now wraps onto 2 lines (same row height as the other 2-line rows)
$ws.Rows.Item(15).RowHeight = 33.6

# Update the view: scroll so row 9 is the top visible row, and move the
# active selection down to D18
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D18").Select()
